$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 397.5
$ws.Range("I31").Value = 397.5
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1192.5
$ws.Range("L31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -962.5
$ws.Range("H86").Value = 26318288
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 26318288
$ws.Range("K86").Value = 0
$ws.Range("L86").ClearContents()
$ws.Range("M86").Value = 26318288
$ws.Range("N86").Value = -26320534
$ws.Range("H89").Value = 26318288
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 26318288
$ws.Range("K89").Value = 0
$ws.Range("L89").ClearContents()
$ws.Range("M89").Value = 131591440
$ws.Range("N89").Value = -131602672
$ws.Range("H92").Value = 625.375
$ws.Range("I92").Value = 250.42857
$ws.Range("J92").Value = 3250
$ws.Range("K92").Value = 250.42857
$ws.Range("L92").Value = 3250
$ws.Range("M92").Value = 997.57143
$ws.Range("N92").Value = -5746
$ws.Range("H132").Value = 4803.525
$ws.Range("I132").Value = 4935.7856
$ws.Range("K132").Value = 14807.3568
$ws.Range("M132").Value = -12277.3568
$ws.Range("H137").Value = 4850.407
$ws.Range("I137").Value = 5085
$ws.Range("K137").Value = 15255
$ws.Range("M137").Value = -12705
$ws.Range("H138").Value = 4653.7734
$ws.Range("I138").Value = 1561.4
$ws.Range("J138").Value = 5874.4473
$ws.Range("K138").Value = 4684.200000000001
$ws.Range("L138").Value = 17623.3419
$ws.Range("M138").Value = 455.7999999999993
$ws.Range("N138").Value = -27903.3419

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 27552.436
$ws.Range("I2").Value = 33863.16
$ws.Range("J2").Value = 3098.375
$ws.Range("K2").Value = 33863.16
$ws.Range("L2").Value = 3098.375
$ws.Range("M2").Value = -33750.16
$ws.Range("N2").Value = -3324.375
$ws.Range("H32").Value = 3636.6875
$ws.Range("I32").Value = 3501.2979
$ws.Range("K32").Value = 3501.2979
$ws.Range("M32").Value = -3214.2979
$ws.Range("H61").Value = 4445.4287
$ws.Range("I61").Value = 1848.6666
$ws.Range("K61").Value = 1848.6666
$ws.Range("M61").Value = -1636.6666
$ws.Range("H62").Value = 45000
$ws.Range("J62").Value = 45000
$ws.Range("L62").Value = 45000
$ws.Range("N62").Value = -46248
$ws.Range("H63").Value = 4715
$ws.Range("I63").Value = 4715
$ws.Range("K63").Value = 4715
$ws.Range("M63").Value = -4029
$ws.Range("H65").Value = 45000
$ws.Range("J65").Value = 45000
$ws.Range("L65").Value = 135000
$ws.Range("N65").Value = -141240
$ws.Range("H66").Value = 4715
$ws.Range("I66").Value = 4715
$ws.Range("K66").Value = 23575
$ws.Range("M66").Value = -20143
$ws.Range("H76").Value = 1000000
$ws.Range("J76").Value = 1000000
$ws.Range("L76").Value = 1000000
$ws.Range("N76").Value = -1000676
$ws.Range("H79").Value = 1000000
$ws.Range("J79").Value = 1000000
$ws.Range("L79").Value = 1000000
$ws.Range("N79").Value = -1002340
$ws.Range("H97").Value = 2161.4285
$ws.Range("I97").Value = 2147.889
$ws.Range("K97").Value = 2147.889
$ws.Range("M97").Value = -1651.889
$ws.Range("H102").Value = 2087.1428
$ws.Range("I102").Value = 1518.75
$ws.Range("J102").Value = 5497.5
$ws.Range("K102").Value = 1518.75
$ws.Range("L102").Value = 5497.5
$ws.Range("M102").Value = 103.25
$ws.Range("N102").Value = -8741.5
$ws.Range("H110").Value = 157514.06
$ws.Range("I110").Value = 186523.27
$ws.Range("K110").Value = 186523.27
$ws.Range("M110").Value = -184478.27
$ws.Range("H116").Value = 27552.436
$ws.Range("I116").Value = 33863.16
$ws.Range("J116").Value = 3098.375
$ws.Range("K116").Value = 33863.16
$ws.Range("L116").Value = 3098.375
$ws.Range("M116").Value = -31569.16
$ws.Range("N116").Value = -7686.375
$ws.Range("H132").Value = 5596.077
$ws.Range("I132").Value = 2723.3215
$ws.Range("K132").Value = 8169.9645
$ws.Range("M132").Value = -5639.9645
$ws.Range("H133").Value = 2537500
$ws.Range("J133").Value = 2537500
$ws.Range("L133").Value = 2537500
$ws.Range("N133").Value = -2542560
$ws.Range("H134").Value = 50000
$ws.Range("J134").Value = 50000
$ws.Range("L134").Value = 50000
$ws.Range("N134").Value = -60140
$ws.Range("H135").Value = 90000
$ws.Range("J135").Value = 90000
$ws.Range("L135").Value = 90000
$ws.Range("N135").Value = -100140
$ws.Range("H136").Value = 4445.4287
$ws.Range("I136").Value = 1848.6666
$ws.Range("K136").Value = 5545.9998
$ws.Range("M136").Value = -2995.9998
$ws.Range("H137").Value = 50000
$ws.Range("J137").Value = 50000
$ws.Range("L137").Value = 50000
$ws.Range("N137").Value = -60200
$ws.Range("H138").Value = 29500
$ws.Range("J138").Value = 29500
$ws.Range("L138").Value = 29500
$ws.Range("N138").Value = -39780
$ws.Range("H140").Value = 49980
$ws.Range("J140").Value = 49980
$ws.Range("L140").Value = 49980
$ws.Range("N140").Value = -60340
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").ClearContents()
$ws.Range("N141").Value = 0

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 27552.436
$ws.Range("I3").Value = 33863.16
$ws.Range("J3").Value = 3098.375
$ws.Range("K3").Value = 33863.16
$ws.Range("L3").Value = 3098.375
$ws.Range("M3").Value = -33749.16
$ws.Range("N3").Value = -3326.375
$ws.Range("H86").Value = 3295
$ws.Range("I86").Value = 2491.6667
$ws.Range("J86").Value = 4500
$ws.Range("K86").Value = 2491.6667
$ws.Range("L86").Value = 4500
$ws.Range("M86").Value = -1368.6667
$ws.Range("N86").Value = -6746
$ws.Range("H89").Value = 3295
$ws.Range("I89").Value = 2491.6667
$ws.Range("J89").Value = 4500
$ws.Range("K89").Value = 12458.3335
$ws.Range("L89").Value = 22500
$ws.Range("M89").Value = -6842.333500000001
$ws.Range("N89").Value = -33732
$ws.Range("H105").Value = 39097.332
$ws.Range("I105").Value = 78685
$ws.Range("K105").Value = 78685
$ws.Range("M105").Value = -76938

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3715.3225
$ws.Range("I132").Value = 2906.7727
$ws.Range("J132").Value = 5691.778
$ws.Range("K132").Value = 8720.3181
$ws.Range("L132").Value = 17075.334
$ws.Range("M132").Value = -6190.3181
$ws.Range("N132").Value = -22135.334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 352.25
$ws.Range("I103").Value = 265.75
$ws.Range("J103").Value = 438.75
$ws.Range("K103").Value = 797.25
$ws.Range("L103").Value = 1316.25
$ws.Range("M103").Value = 81.75
$ws.Range("N103").Value = -3074.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 59852
$ws.Range("J95").Value = 59852
$ws.Range("L95").Value = 59852
$ws.Range("N95").Value = -65344
$ws.Range("H97").Value = 4600.5186
$ws.Range("I97").Value = 5446.9546
$ws.Range("J97").Value = 876.2
$ws.Range("K97").Value = 5446.9546
$ws.Range("L97").Value = 876.2
$ws.Range("M97").Value = -4950.9546
$ws.Range("N97").Value = -1868.2
$ws.Range("H99").Value = 3939.5
$ws.Range("I99").Value = 3939.5
$ws.Range("K99").Value = 3939.5
$ws.Range("M99").Value = -1693.5
$ws.Range("H111").Value = 35992.285
$ws.Range("J111").Value = 35992.285
$ws.Range("L111").Value = 35992.285
$ws.Range("N111").Value = -42126.285
$ws.Range("H122").Value = 132158.62
$ws.Range("I122").Value = 207703.8
$ws.Range("J122").Value = 6250
$ws.Range("K122").Value = 623111.3999999999
$ws.Range("L122").Value = 18750
$ws.Range("M122").Value = -620661.3999999999
$ws.Range("N122").Value = -23650
$ws.Range("H132").Value = 266092.88
$ws.Range("I132").Value = 296697.94
$ws.Range("J132").Value = 5949.75
$ws.Range("K132").Value = 890093.8200000001
$ws.Range("L132").Value = 17849.25
$ws.Range("M132").Value = -887563.8200000001
$ws.Range("N132").Value = -22909.25
$ws.Range("H133").Value = 50000
$ws.Range("J133").Value = 50000
$ws.Range("L133").Value = 50000
$ws.Range("N133").Value = -60120

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 771692.4
$ws.Range("I40").Value = 910681.9399999999
$ws.Range("K40").Value = 910681.9399999999
$ws.Range("M40").Value = -910545.9399999999
$ws.Range("H55").Value = 516.8372000000001
$ws.Range("I55").Value = 437.0606
$ws.Range("J55").Value = 780.1
$ws.Range("K55").Value = 437.0606
$ws.Range("L55").Value = 780.1
$ws.Range("M55").Value = -264.0606
$ws.Range("N55").Value = -1126.1
$ws.Range("H61").Value = 4265.121
$ws.Range("I61").Value = 3128.95
$ws.Range("K61").Value = 3128.95
$ws.Range("M61").Value = -2926.95
$ws.Range("H110").Value = 41732.75
$ws.Range("J110").Value = 41732.75
$ws.Range("L110").Value = 41732.75
$ws.Range("N110").Value = -49912.75
$ws.Range("H113").Value = 4265.121
$ws.Range("I113").Value = 3128.95
$ws.Range("K113").Value = 3128.95
$ws.Range("M113").Value = -958.9499999999998
$ws.Range("H132").Value = 7543.8887
$ws.Range("I132").Value = 4615
$ws.Range("K132").Value = 13845
$ws.Range("M132").Value = -11315
$ws.Range("H136").Value = 4310.7334
$ws.Range("I136").Value = 3605.9092
$ws.Range("J136").Value = 6249
$ws.Range("K136").Value = 10817.7276
$ws.Range("L136").Value = 18747
$ws.Range("M136").Value = -8267.7276
$ws.Range("N136").Value = -23847

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 15875996
$ws.Range("I122").Value = 20835490
$ws.Range("J122").Value = 5612.2666
$ws.Range("K122").Value = 62506470
$ws.Range("L122").Value = 16836.7998
$ws.Range("M122").Value = -62504020
$ws.Range("N122").Value = -21736.7998
